$wb = $excel.ActiveWorkbook

# --- Update the "Hoja1" note text in cell A1 ---
$wsNota = $wb.Worksheets.Item("Hoja1")
$cell = $wsNota.Range("A1")
$old = $cell.Value()
$oldLine1 = "✅ 1000 Bs = 15.16 = 62956.95 pesos"
$oldLine2 = "✅ 62956.95 pesos = 15.0 = 962.2 Bs"
$newLine1 = "✅ 1000 Bs = 14.95 = 62152.47 pesos"
$newLine2 = "✅ 62152.47 pesos = 14.94 = 980.77 Bs"
$new = $old.Replace($oldLine1, $newLine1).Replace($oldLine2, $newLine2)
$cell.Value = $new

# --- Update the rate figures on the "tasas" sheet ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 66.90000000000001
$wsTasas.Range("O10").Value = 4158
$wsTasas.Range("N12").Value = 4160
$wsTasas.Range("O12").Value = 65.645
